$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Stage the distinct formats we will need, by copying them off the
#    existing (pre-edit) cells that already carry them, onto scratch cells
#    far outside the used range. We read these BEFORE any destructive edit.
# ---------------------------------------------------------------------------
$stage = @{
  "s2" = "BZ101"   # row10 header: bold white/black, border1, wrap, center/center
  "s3" = "BZ102"   # row10 header no-wrap: border1, center/center
  "s5" = "BZ103"   # row10 header border0: center/center
  "s7" = "BZ104"   # row9 title: bold black on white, border0, h-center
  "s8" = "BZ105"   # row9 title: bold black no-fill, border0, h-center
  "s9" = "BZ106"   # row9 title: bold black no-fill, border4, h-center
  "s4" = "BZ107"   # border3, center/center
  "s6" = "BZ108"   # border2, center/center (new alignment combo)
}

$ws.Range("A10").Copy()
$ws.Range($stage["s2"]).PasteSpecial(-4122) | Out-Null

$ws.Range("B10").Copy()
$ws.Range($stage["s3"]).PasteSpecial(-4122) | Out-Null

$ws.Range("C10").Copy()
$ws.Range($stage["s5"]).PasteSpecial(-4122) | Out-Null

$ws.Range("A9").Copy()
$ws.Range($stage["s7"]).PasteSpecial(-4122) | Out-Null

$ws.Range("E9").Copy()
$ws.Range($stage["s8"]).PasteSpecial(-4122) | Out-Null

$ws.Range("L9").Copy()
$ws.Range($stage["s9"]).PasteSpecial(-4122) | Out-Null

$ws.Range("AL10").Copy()
$ws.Range($stage["s4"]).PasteSpecial(-4122) | Out-Null

$ws.Range("AK10").Copy()
$ws.Range($stage["s6"]).PasteSpecial(-4122) | Out-Null
$ws.Range($stage["s6"]).HorizontalAlignment = -4108

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Unmerge the old row-9 merges, then clear the old header rows entirely.
# ---------------------------------------------------------------------------
$ws.Range("A9:AM9").UnMerge()
$ws.Range("A9:AM10").Clear() | Out-Null

function Stamp($addr, $styleKey) {
  $ws.Range($stage[$styleKey]).Copy()
  $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. Rebuild row 9 (group titles) and row 10 (sub-headers).
#    Cell *values* are written in the exact left-to-right / top-to-bottom
#    order of the target workbook so newly-introduced shared strings land
#    at the same indices as the authored edit.
# ---------------------------------------------------------------------------

# --- Group 1: A9:F9 "ORDEN DE AUDITORÍA" ---
Stamp "A9" "s7"; Stamp "B9" "s7"; Stamp "C9" "s7"; Stamp "D9" "s7"; Stamp "E9" "s7"; Stamp "F9" "s7"
$ws.Range("A9").Value = "ORDEN DE AUDITORÍA"

Stamp "A10" "s2"
$ws.Range("A10").Value = "OFICIO ORDEN"
Stamp "B10" "s3"
$ws.Range("B10").Value = "FECHA OFICIO"
Stamp "C10" "s5"
$ws.Range("C10").Value = "OFICIO NOTIFICACIÓN"
Stamp "D10" "s5"
$ws.Range("D10").Value = "FECHA VENCIMIENTO"
Stamp "E10" "s5"
$ws.Range("E10").Value = "OFICIO SOLICITUD"
Stamp "F10" "s5"
$ws.Range("F10").Value = "FECHA SOLICITUD"

$ws.Range("A9:F9").Merge() | Out-Null

# --- Group 2: H9:M9 "NOTIFICACIÓN A ÁREAS" ---
Stamp "H9" "s8"; Stamp "I9" "s8"; Stamp "J9" "s8"; Stamp "K9" "s8"; Stamp "L9" "s8"; Stamp "M9" "s8"
$ws.Range("H9").Value = "NOTIFICACIÓN A ÁREAS"

Stamp "H10" "s2"
$ws.Range("H10").Value = "UNIDAD ADMINISTRATIVA"
Stamp "I10" "s2"
$ws.Range("I10").Value = "OFICIO"
Stamp "J10" "s2"
$ws.Range("J10").Value = "FECHA OFICIO"
Stamp "K10" "s2"
$ws.Range("K10").Value = "FECHA RECIBIDO"
Stamp "L10" "s2"
$ws.Range("L10").Value = "FECHA VENCIMIENTO"
Stamp "M10" "s2"
$ws.Range("M10").Value = "PRÓRROGA"

$ws.Range("H9:M9").Merge() | Out-Null

# --- Group 3: O9:T9 "ÁREA NOTIFICADA" ---
Stamp "O9" "s9"; Stamp "P9" "s9"; Stamp "Q9" "s9"; Stamp "R9" "s9"; Stamp "S9" "s9"; Stamp "T9" "s9"
$ws.Range("O9").Value = "ÁREA NOTIFICADA"

Stamp "O10" "s2"
$ws.Range("O10").Value = "UNIDAD ADMINISTRATIVA"
Stamp "P10" "s2"
$ws.Range("P10").Value = "OFICIO"
Stamp "Q10" "s2"
$ws.Range("Q10").Value = "FECHA OFICIO"
Stamp "R10" "s2"
$ws.Range("R10").Value = "FECHA RECIBIDO"
Stamp "S10" "s2"
$ws.Range("S10").Value = "FECHA VENCIMIENTO"
Stamp "T10" "s2"
$ws.Range("T10").Value = "PRÓRROGA"

$ws.Range("O9:T9").Merge() | Out-Null

# --- Group 4: V9:AB9 "CONTESTACIÓN ÓRGANO FISCALIZADOR" ---
Stamp "V9" "s9"; Stamp "W9" "s9"; Stamp "X9" "s9"; Stamp "Y9" "s9"; Stamp "Z9" "s9"; Stamp "AA9" "s9"; Stamp "AB9" "s9"
$ws.Range("V9").Value = "CONTESTACIÓN ÓRGANO FISCALIZADOR"

Stamp "V10" "s2"
$ws.Range("V10").Value = "ÓRGANO"
Stamp "W10" "s2"
$ws.Range("W10").Value = "OFICIO"
Stamp "X10" "s2"
$ws.Range("X10").Value = "FOLIO SIGA"
Stamp "Y10" "s2"
$ws.Range("Y10").Value = "ENTREGA"
Stamp "Z10" "s2"
$ws.Range("Z10").Value = "FECHA OFICIO"
Stamp "AA10" "s2"
$ws.Range("AA10").Value = "FECHA RECIBIDO"
Stamp "AB10" "s2"
$ws.Range("AB10").Value = "FECHA VENCIMENTO"

$ws.Range("V9:AB9").Merge() | Out-Null

# --- New row 23 marker cell (introduces the lone-space shared string) ---
$ws.Range("AJ23").Value = " "

# --- Group 5 (NEW): AD9:AJ9 "RESPUESTA A CONTESTACIÓN ÓRGANO FISCALIZADOR" ---
# Sub-header values are written before the group title so the new shared
# strings land in the same order as the authored edit.
Stamp "AD9" "s9"; Stamp "AE9" "s9"; Stamp "AF9" "s9"; Stamp "AG9" "s9"; Stamp "AH9" "s9"; Stamp "AI9" "s9"; Stamp "AJ9" "s9"

Stamp "AD10" "s2"
$ws.Range("AD10").Value = "ÓRGANO ORIGEN"
Stamp "AE10" "s2"
$ws.Range("AE10").Value = "ÓRGANO DESTINO"
Stamp "AF10" "s2"
$ws.Range("AF10").Value = "OFICIO"
Stamp "AG10" "s2"
$ws.Range("AG10").Value = "FOLIO SIGA"
Stamp "AH10" "s2"
$ws.Range("AH10").Value = "FECHA OFICIO"
Stamp "AI10" "s2"
$ws.Range("AI10").Value = "FECHA RECIBIDO"
Stamp "AJ10" "s2"
$ws.Range("AJ10").Value = "FECHA VENCIMENTO"

$ws.Range("AD9").Value = "RESPUESTA A CONTESTACIÓN ÓRGANO FISCALIZADOR"
$ws.Range("AD9:AJ9").Merge() | Out-Null

# --- Group 6: AL9:AR9 "ADMINISTRACIÓN DE RESULTADOS" ---
Stamp "AL9" "s8"; Stamp "AM9" "s8"; Stamp "AN9" "s8"; Stamp "AO9" "s8"; Stamp "AP9" "s8"; Stamp "AQ9" "s8"; Stamp "AR9" "s8"
$ws.Range("AL9").Value = "ADMINISTRACIÓN DE RESULTADOS"

Stamp "AL10" "s2"
$ws.Range("AL10").Value = "TIPO DE RESULTADO"
Stamp "AM10" "s2"
$ws.Range("AM10").Value = "ESTATUS DE LOS RESULTADOS"
Stamp "AN10" "s2"
$ws.Range("AN10").Value = "MONTO"
Stamp "AO10" "s2"
$ws.Range("AO10").Value = "CLAVE DE RESULTADO"
Stamp "AP10" "s2"
$ws.Range("AP10").Value = "RESULTADO/OBSERVACIÓN"
Stamp "AQ10" "s2"
$ws.Range("AQ10").Value = "RESULTADO SUPERVIVIENTE"
Stamp "AR10" "s2"
$ws.Range("AR10").Value = "NUMERO RESULTADO"

$ws.Range("AL9:AR9").Merge() | Out-Null

# --- Group 7: AT9:AY9 "SE SOLICITA INFORMACIÓN PRELIMINAR A LAS DIRECCIONES DE LA SFYTGE:" ---
Stamp "AT9" "s7"; Stamp "AU9" "s7"; Stamp "AV9" "s7"; Stamp "AW9" "s7"; Stamp "AX9" "s7"; Stamp "AY9" "s7"
$ws.Range("AT9").Value = "SE SOLICITA INFORMACIÓN PRELIMINAR A LAS DIRECCIONES DE LA SFYTGE:"

Stamp "AT10" "s2"
$ws.Range("AT10").Value = "ÁREA NOTIFICADA"
Stamp "AU10" "s3"
$ws.Range("AU10").Value = "No. DE OFICIO"
Stamp "AV10" "s3"
$ws.Range("AV10").Value = "FECHA"
Stamp "AW10" "s6"
$ws.Range("AW10").Value = "SOLICITUD"
Stamp "AX10" "s4"
$ws.Range("AX10").Value = "No. DE OFICIO"
Stamp "AY10" "s4"
$ws.Range("AY10").Value = "FECHA RECIBIDO"

$ws.Range("AT9:AY9").Merge() | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Clean up the staging cells now that every format has been consumed.
# ---------------------------------------------------------------------------
$ws.Range("BZ101:BZ108").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 5. Row height for the rebuilt header row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 15.5

# ---------------------------------------------------------------------------
# 6. Column widths (best effort; engine quantizes to pixel grid).
# ---------------------------------------------------------------------------
$colWidths = @{
  1=19.1796875; 2=22.1796875; 3=22.1796875; 4=22.1796875; 5=22.1796875; 6=22.1796875;
  8=28.90625; 9=15; 10=12.54296875; 11=14.7265625; 12=18.81640625; 13=10.36328125;
  15=28.36328125; 16=20.1796875; 17=16.1796875; 18=19.36328125; 19=20.90625; 20=20.90625;
  22=22; 23=13.90625; 24=11.453125; 25=11.453125; 26=12.54296875; 27=14.7265625; 28=18.26953125;
  30=22.54296875; 31=19.1796875; 32=15.90625; 33=12.453125; 34=19.453125; 35=19.453125; 36=19.453125;
  38=18; 39=25.81640625; 40=7.54296875; 41=21.1796875; 42=32.81640625; 43=24.453125; 44=19;
  46=27.54296875; 47=15; 48=12.1796875; 49=18.81640625; 50=15.08984375; 51=14.7265625
}
foreach ($c in $colWidths.Keys) {
  $ws.Columns.Item($c).ColumnWidth = $colWidths[$c]
}

# ---------------------------------------------------------------------------
# 7. Sheet view: drop frozen/topLeftCell scroll, select B2.
# ---------------------------------------------------------------------------
$ws.Range("B2").Select()
